# Fix general intent submodules:
#  - Update the "BANK_ACC has AMOUNT left." response text on the
#    GENERAL_INTENTS sheet to "BANK_ACC bank account has AMOUNT left."
#  - Switch the saved selection/active tab from BANK_BALANCE back to
#    GENERAL_INTENTS (with B11 selected there).

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("GENERAL_INTENTS")

# Update the response text in B5 of GENERAL_INTENTS.
$wsGeneral.Range("B5").Value = "BANK_ACC bank account has AMOUNT left."

# Make GENERAL_INTENTS the active/selected sheet again (it had lost that
# status to BANK_BALANCE) and restore its selection to cell B11.
$wsGeneral.Activate()
$wsGeneral.Range("B11").Select()

# BANK_BALANCE keeps its existing selection (E18); it simply stops being
# the tab that is active/selected once GENERAL_INTENTS is activated above.
